$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 652.3935292841447
$ws.Range("E2").Value = 289624.7207921735
$ws.Range("I2").Value = 264280.8358355595
$ws.Range("L2").Value = 269030.948845862
$ws.Range("M2").Value = 116987.63646295
$ws.Range("N2").Value = 71604.66739785175
$ws.Range("O2").Value = 69620.67602513026

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 1448.958715432552
$ws.Range("B2").Value = 35880.83082778667
$ws.Range("E2").Value = 164812.2831047642
$ws.Range("I2").Value = 209218.95419975
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 61314.03182440259
$ws.Range("N2").Value = 19801.59487199476
$ws.Range("O2").Value = 10580.29689353687

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 49430.84256000001
$ws.Range("B2").Value = 29803.23071747997
$ws.Range("E2").Value = 56470.04597280241
$ws.Range("I2").Value = 153491.5861848302
$ws.Range("M2").Value = 59790.19096888593
$ws.Range("N2").Value = 22966.81720506174
$ws.Range("O2").Value = 57371.12014516797
